# Legenda.xlsx update:
#  - On sheet "E7X", the fund "Name" column (B2:B9) is simplified by
#    dropping the "Eurizon Fund - " / "Eurizon Investment Sicav - "
#    prefix, keeping only the strategy name.
#  - The active sheet/selection moves from "Scenari" to "E7X" (cell B10).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("E7X")

$ws.Range("B2").Value = "Dynamic Asset Allocation"
$ws.Range("B3").Value = "Flexible Multistrategy"
$ws.Range("B4").Value = "Active Allocation"
$ws.Range("B5").Value = "Strategic Allocation"
$ws.Range("B6").Value = "Sustainable Multiasset"
$ws.Range("B7").Value = "Conservative Allocation"
$ws.Range("B8").Value = "Absolute Prudent"
$ws.Range("B9").Value = "Flexible Macro Factors"

# Make E7X the active sheet with B10 selected (matches the new
# workbookView/sheetView state captured in the saved file).
$ws.Activate()
$ws.Range("B10").Select()
